$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-11-09 Sunday" "2025-11-10 Monday"

Replace-Text "112×7=" "957×5="
Replace-Text "151×9=" "752×8="
Replace-Text "153×7=" "559×2="
Replace-Text "474×7=" "114×4="
Replace-Text "784×7=" "845×6="

Replace-Text "202×2=" "582×5="
Replace-Text "105×6=" "172×2="
Replace-Text "525×6=" "967×7="
Replace-Text "302×5=" "195×2="
Replace-Text "979×8=" "367×4="

Replace-Text "562×9=" "287×8="
Replace-Text "311×6=" "194×4="
Replace-Text "329×9=" "939×9="
Replace-Text "720×5=" "363×7="
Replace-Text "622×9=" "324×9="

Replace-Text "755×6=" "136×2="
Replace-Text "613×4=" "800×8="
Replace-Text "580×8=" "944×4="
Replace-Text "163×6=" "629×2="
Replace-Text "166×3=" "199×4="

Replace-Text "627×8=" "259×6="
Replace-Text "314×6=" "506×2="
Replace-Text "594×3=" "335×2="
Replace-Text "455×7=" "575×9="
Replace-Text "528×7=" "653×7="
